# EPBDS-11024: add more tests for SmartLookup table should not be transposed during compilation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Change the header text of the existing table (row 6) ---
$ws.Range("C6").Value = "SmartLookup String test(  String param1,  String  param2, String param3, int param4,String  param5)"

# --- New table #1 (rows 16-20): SmartLookup  Integer mySmartLookup(...) ---
$ws.Cells.Item(16, 3).Value = "SmartLookup  Integer mySmartLookup(String  carNum, CarValue carValue , Double param3, Double param4)"
$ws.Range("C16:F16").Merge()
$ws.Range("C16:F16").HorizontalAlignment = -4131

$ws.Cells.Item(17, 3).Value = "num"
$ws.Range("C17:C18").Merge()
$ws.Range("C17:C18").HorizontalAlignment = -4108

$ws.Cells.Item(17, 4).Value = "test"
$ws.Range("D17:E17").Merge()
$ws.Range("D17:E17").HorizontalAlignment = -4108

$ws.Cells.Item(17, 6).Value = "pam"
$ws.Cells.Item(17, 6).HorizontalAlignment = -4108

$ws.Cells.Item(18, 4).Value = "50.0"
$ws.Cells.Item(18, 5).Value = "60.0"
$ws.Cells.Item(18, 6).Value = 7

$ws.Cells.Item(19, 3).Value = 3
$ws.Cells.Item(19, 4).Value = 2
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = 3

$ws.Cells.Item(20, 3).Value = 45
$ws.Cells.Item(20, 4).Value = 2
$ws.Cells.Item(20, 5).Value = 1
$ws.Cells.Item(20, 6).Value = 4

$ws.Range("C16:F20").Borders.Item(1).LineStyle = 1
$ws.Range("C16:F20").Borders.Item(1).Weight = 2
$ws.Range("C16:F20").Borders.Item(2).LineStyle = 1
$ws.Range("C16:F20").Borders.Item(2).Weight = 2
$ws.Range("C16:F20").Borders.Item(3).LineStyle = 1
$ws.Range("C16:F20").Borders.Item(3).Weight = 2
$ws.Range("C16:F20").Borders.Item(4).LineStyle = 1
$ws.Range("C16:F20").Borders.Item(4).Weight = 2
$ws.Range("C16:F20").Borders.Item(7).LineStyle = 1
$ws.Range("C16:F20").Borders.Item(7).Weight = 2
$ws.Range("C16:F20").Borders.Item(8).LineStyle = 1
$ws.Range("C16:F20").Borders.Item(8).Weight = 2
$ws.Range("C16:F20").Borders.Item(9).LineStyle = 1
$ws.Range("C16:F20").Borders.Item(9).Weight = 2
$ws.Range("C16:F20").Borders.Item(10).LineStyle = 1
$ws.Range("C16:F20").Borders.Item(10).Weight = 2
$ws.Range("C16:F20").Borders.Item(11).LineStyle = 1
$ws.Range("C16:F20").Borders.Item(11).Weight = 2
$ws.Range("C16:F20").Borders.Item(12).LineStyle = 1
$ws.Range("C16:F20").Borders.Item(12).Weight = 2

# --- Explanatory text block (rows 23-26) ---
$ws.Cells.Item(23, 3).Value = "Datatype CarValue <String>"
$ws.Cells.Item(24, 3).Value = "num"
$ws.Cells.Item(25, 3).Value = "param"
$ws.Cells.Item(26, 3).Value = "pam"

# --- New table #2 (rows 28-32): Transposed version SmartLookup  Integer mySmartLookupTrans(...) ---
$ws.Cells.Item(28, 3).Value = "SmartLookup  Integer mySmartLookupTrans(String  carNum, CarValue carValue , Double param3, Double param4)"
$ws.Range("C28:F28").Merge()
$ws.Range("C28:F28").HorizontalAlignment = -4131

$ws.Cells.Item(29, 3).Value = "num"
$ws.Range("C29:D29").Merge()
$ws.Range("C29:D29").HorizontalAlignment = -4108
$ws.Cells.Item(29, 5).Value = 3
$ws.Cells.Item(29, 6).Value = 45

$ws.Cells.Item(30, 3).Value = "test"
$ws.Range("C30:C31").Merge()
$ws.Range("C30:C31").HorizontalAlignment = -4108
$ws.Cells.Item(30, 4).Value = "50.0"
$ws.Cells.Item(30, 5).Value = 2
$ws.Cells.Item(30, 6).Value = 2

$ws.Cells.Item(31, 4).Value = "60.0"
$ws.Cells.Item(31, 5).Value = 3
$ws.Cells.Item(31, 6).Value = 1

$ws.Cells.Item(32, 3).Value = "pam"
$ws.Cells.Item(32, 4).Value = 7
$ws.Cells.Item(32, 5).Value = 3
$ws.Cells.Item(32, 6).Value = 4

$ws.Range("C28:F32").Borders.Item(1).LineStyle = 1
$ws.Range("C28:F32").Borders.Item(1).Weight = 2
$ws.Range("C28:F32").Borders.Item(2).LineStyle = 1
$ws.Range("C28:F32").Borders.Item(2).Weight = 2
$ws.Range("C28:F32").Borders.Item(3).LineStyle = 1
$ws.Range("C28:F32").Borders.Item(3).Weight = 2
$ws.Range("C28:F32").Borders.Item(4).LineStyle = 1
$ws.Range("C28:F32").Borders.Item(4).Weight = 2
$ws.Range("C28:F32").Borders.Item(7).LineStyle = 1
$ws.Range("C28:F32").Borders.Item(7).Weight = 2
$ws.Range("C28:F32").Borders.Item(8).LineStyle = 1
$ws.Range("C28:F32").Borders.Item(8).Weight = 2
$ws.Range("C28:F32").Borders.Item(9).LineStyle = 1
$ws.Range("C28:F32").Borders.Item(9).Weight = 2
$ws.Range("C28:F32").Borders.Item(10).LineStyle = 1
$ws.Range("C28:F32").Borders.Item(10).Weight = 2
$ws.Range("C28:F32").Borders.Item(11).LineStyle = 1
$ws.Range("C28:F32").Borders.Item(11).Weight = 2
$ws.Range("C28:F32").Borders.Item(12).LineStyle = 1
$ws.Range("C28:F32").Borders.Item(12).Weight = 2

# --- View state: selection on the new second table header ---
$ws.Range("C28:F28").Select()

Write-Host "EPBDS-11024 transposed-table tests added"
